$wb = $excel.ActiveWorkbook

# Row 55 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 1563206
$ws.Range("I55").Value = 1883.3334
$ws.Range("J55").Value = 2232344.2
$ws.Range("K55").Value = 1883.3334
$ws.Range("L55").Value = 2232344.2
$ws.Range("M55").Value = -1669.3334
$ws.Range("N55").Value = -2232772.2

# Row 118 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 250
$ws.Range("I118").Value = 250
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 750
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 907
$ws.Range("N118").ClearContents()

# Row 141 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3289.7354
$ws.Range("I141").Value = 1806.24
$ws.Range("J141").Value = 7410.5557
$ws.Range("K141").Value = 5418.72
$ws.Range("L141").Value = 22231.6671
$ws.Range("M141").Value = -238.7200000000003
$ws.Range("N141").Value = -32591.6671

# Row 32 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 629.11
$ws.Range("I32").Value = 623.5851
$ws.Range("J32").Value = 715.6667
$ws.Range("K32").Value = 623.5851
$ws.Range("L32").Value = 715.6667
$ws.Range("M32").Value = -336.5851
$ws.Range("N32").Value = -1289.6667

# Row 74 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 742.629
$ws.Range("I74").Value = 799.3269
$ws.Range("J74").Value = 447.8
$ws.Range("K74").Value = 799.3269
$ws.Range("L74").Value = 447.8
$ws.Range("M74").Value = 74.67309999999998
$ws.Range("N74").Value = -2195.8

# Row 77 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 742.629
$ws.Range("I77").Value = 799.3269
$ws.Range("J77").Value = 447.8
$ws.Range("K77").Value = 3996.6345
$ws.Range("L77").Value = 2239
$ws.Range("M77").Value = 371.3654999999999
$ws.Range("N77").Value = -10975

# Row 101 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

# Row 102 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1677.75
$ws.Range("I102").Value = 1052.5
$ws.Range("K102").Value = 1052.5
$ws.Range("M102").Value = 569.5

# Row 132 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1636017.8
$ws.Range("I132").Value = 1670.4333
$ws.Range("K132").Value = 5011.2999
$ws.Range("M132").Value = -2481.2999

# Row 134 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2587316
$ws.Range("I134").Value = 963.0606
$ws.Range("K134").Value = 2889.1818
$ws.Range("M134").Value = -354.1818000000003

# Row 97 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 30000
$ws.Range("J97").Value = 30000
$ws.Range("L97").Value = 30000
$ws.Range("N97").Value = -31982

# Row 105 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3594.2856
$ws.Range("I105").Value = 3617.7778
$ws.Range("J105").Value = 3256
$ws.Range("K105").Value = 3617.7778
$ws.Range("L105").Value = 3256
$ws.Range("M105").Value = -1870.7778
$ws.Range("N105").Value = -6750

# Row 122 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 19232262
$ws.Range("I122").Value = 31251288
$ws.Range("J122").Value = 1820
$ws.Range("K122").Value = 93753864
$ws.Range("L122").Value = 5460
$ws.Range("M122").Value = -93751414
$ws.Range("N122").Value = -10360

# Row 34 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 566.8
$ws.Range("I34").Value = 600.5
$ws.Range("J34").Value = 554.5454999999999
$ws.Range("K34").Value = 1801.5
$ws.Range("L34").Value = 1663.6365
$ws.Range("M34").Value = -1717.5
$ws.Range("N34").Value = -1831.6365

# Row 39 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 200
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

# Row 55 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 599.6667
$ws.Range("J55").Value = 1500
$ws.Range("L55").Value = 4500
$ws.Range("N55").Value = -4854

# Row 116 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 111112690
$ws.Range("I116").Value = 457.33334
$ws.Range("J116").Value = 166668800
$ws.Range("K116").Value = 1372.00002
$ws.Range("L116").Value = 500006400
$ws.Range("M116").Value = 2069.99998
$ws.Range("N116").Value = -500013284

# Row 118 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 29102.75
$ws.Range("I118").Value = 1500
$ws.Range("J118").Value = 38303.668
$ws.Range("K118").Value = 4500
$ws.Range("L118").Value = 114911.004
$ws.Range("M118").Value = -3257
$ws.Range("N118").Value = -117397.004

# Row 121 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 9259934
$ws.Range("I121").Value = 500
$ws.Range("J121").Value = 12346412
$ws.Range("K121").Value = 1500
$ws.Range("L121").Value = 37039236
$ws.Range("M121").Value = -190
$ws.Range("N121").Value = -37041856

# Row 122 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 10421010
$ws.Range("I122").Value = 48077240
$ws.Range("J122").Value = 5457.213
$ws.Range("K122").Value = 432695160
$ws.Range("L122").Value = 49114.917
$ws.Range("M122").Value = -432692710
$ws.Range("N122").Value = -54014.917

# Row 131 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 776.949
$ws.Range("I131").Value = 371.5
$ws.Range("J131").Value = 833.52325
$ws.Range("K131").Value = 1114.5
$ws.Range("L131").Value = 2500.56975
$ws.Range("M131").Value = 3925.5
$ws.Range("N131").Value = -12580.56975

# Row 86 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 14800
$ws.Range("J86").Value = 14800
$ws.Range("L86").Value = 14800
$ws.Range("N86").Value = -17172

# Row 89 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H89").Value = 14800
$ws.Range("J89").Value = 14800
$ws.Range("L89").Value = 44400
$ws.Range("N89").Value = -56256

# Row 126 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2334.5715
$ws.Range("I126").Value = 1642.4
$ws.Range("J126").Value = 2550.875
$ws.Range("K126").Value = 4927.200000000001
$ws.Range("L126").Value = 7652.625
$ws.Range("M126").Value = -2457.200000000001
$ws.Range("N126").Value = -12592.625

# Row 40 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6580770.5
$ws.Range("I40").Value = 1372.2
$ws.Range("J40").Value = 31253512
$ws.Range("K40").Value = 1372.2
$ws.Range("L40").Value = 31253512
$ws.Range("M40").Value = -1236.2
$ws.Range("N40").Value = -31253784

# Row 88 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 13900
$ws.Range("I88").Value = 13900
$ws.Range("K88").Value = 13900
$ws.Range("M88").Value = -13472

# Row 91 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H91").Value = 13900
$ws.Range("I91").Value = 13900
$ws.Range("K91").Value = 13900
$ws.Range("M91").Value = -12418

# Row 132 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 15389485
$ws.Range("I132").Value = 27779278
$ws.Range("K132").Value = 83337834
$ws.Range("M132").Value = -83335304

# Row 56 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 8428
$ws.Range("I56").Value = 8428
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 8428
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -7714
$ws.Range("N56").ClearContents()

# Row 113 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 52632148
$ws.Range("I113").Value = 66667144
$ws.Range("J113").Value = 911.25
$ws.Range("K113").Value = 200001432
$ws.Range("L113").Value = 2733.75
$ws.Range("M113").Value = -199999262
$ws.Range("N113").Value = -7073.75

# Row 126 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2362
$ws.Range("I126").Value = 1599.2
$ws.Range("K126").Value = 4797.6
$ws.Range("M126").Value = -2327.6

# Row 132 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 36306.395
$ws.Range("I132").Value = 46448.176
$ws.Range("J132").Value = 12980.3
$ws.Range("K132").Value = 139344.528
$ws.Range("L132").Value = 38940.89999999999
$ws.Range("M132").Value = -136814.528
$ws.Range("N132").Value = -44000.89999999999
